$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.795.41'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '2.102.11'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.93'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.25'
$ws.Range('E7').Value = '  +3.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.389'
$ws.Range('E9').Value = '  +1.92%  '
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('E12').Value = '  +5.94%  '
$ws.Range('D13').Value = '2.412.88'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.95'
$ws.Range('E14').Value = '  -1.69%  '
$ws.Range('E15').Value = '  +2.51%  '
$ws.Range('E16').Value = '  +0.63%  '
$ws.Range('D17').Value = '2.075.42'
$ws.Range('E17').Value = '  -0.07%  '
$ws.Range('D18').Value = '38.760.71'
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.73'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('D21').Value = '0.0₃0847'
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.45'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.50'
$ws.Range('E24').Value = '  +2.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.65'
$ws.Range('E26').Value = '  +2.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.57'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('E29').Value = '  +2.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.34'
$ws.Range('E30').Value = '  +1.39%  '
$ws.Range('E31').Value = '  +9.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.121'
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.58'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('B34').Value = 'THORChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.15'
$ws.Range('E34').Value = '  +13.25%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.76'
$ws.Range('E35').Value = '  -0.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0614'
$ws.Range('E36').Value = '  +1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.37'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.00'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('E41').Value = '  +3.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '101.38'
$ws.Range('E42').Value = '  +1.21%  '
$ws.Range('D43').Value = '1.527.65'
$ws.Range('E44').Value = '  +7.68%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.80'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.80'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0914'
$ws.Range('E47').Value = '  -0.94%  '
$ws.Range('E48').Value = '  +5.35%  '
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('D51').Value = '2.300.76'
$ws.Range('E51').Value = '  +1.12%  '
